# Applies the "coupling.pptx" tutorial-figure fixups:
#  - refresh the cached date placeholder text (slide master + all 11 layouts)
#  - rename "Master" -> "Parent" and "Sub-App" -> "Child-App" on both slides

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "1/13/21") {
                $sh.TextFrame.TextRange.Text = "9/10/23"
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's date placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Rename the diagram boxes on every slide.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -eq "Master") {
                $sh.TextFrame.TextRange.Text = "Parent"
            } elseif ($t -eq "Sub-App") {
                $sh.TextFrame.TextRange.Text = "Child-App"
            }
        }
    }
}
